# Update the "Algo" result values (column B) and a few column A values
# to reflect the updated KNN imputation algorithm results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = 6.922
$ws.Range("B3").Value  = 6.798999999999999
$ws.Range("B5").Value  = 6.834999999999999
$ws.Range("A9").Value  = -21.174
$ws.Range("B11").Value = 6.93
$ws.Range("B12").Value = 6.582000000000001
$ws.Range("A13").Value = -21.832
$ws.Range("A16").Value = -20.798
$ws.Range("A18").Value = -21.868
$ws.Range("A20").Value = -21.664
$ws.Range("B21").Value = 7.007000000000001
